# Natmi following Dr Hou advice
# The sending-cluster set now includes "ECs" (previously only FAPs/sCs were
# used as senders), so the LR-pair table grows from 2 senders x 3 targets (6 rows)
# to 3 senders x 3 targets (9 rows), and every row is recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe out the old 6 data rows (A2:T7) so the sheet can be rebuilt cleanly with
# the new 9-row table (A2:T10).
$ws.Range("A2:T7").ClearContents()

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ncam1"
$ws.Cells.Item(2, 3).Value = "Fgfr1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.9404873333333333
$ws.Cells.Item(2, 8).Value = 2.821462
$ws.Cells.Item(2, 9).Value = 0.02000383747045655
$ws.Cells.Item(2, 10).Value = 0.02000383747045654
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 4.675378666666666
$ws.Cells.Item(2, 14).Value = 14.026136
$ws.Cells.Item(2, 15).Value = 0.03681964474327726
$ws.Cells.Item(2, 16).Value = 0.03681964474327726
$ws.Cells.Item(2, 17).Value = 4.397134414536889
$ws.Cells.Item(2, 18).Value = 39.574209730832
$ws.Cells.Item(2, 19).Value = 0.0007365341891644681
$ws.Cells.Item(2, 20).Value = 0.000736534189164468

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ncam1"
$ws.Cells.Item(3, 3).Value = "Fgfr1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.9404873333333333
$ws.Cells.Item(3, 8).Value = 2.821462
$ws.Cells.Item(3, 9).Value = 0.02000383747045655
$ws.Cells.Item(3, 10).Value = 0.02000383747045654
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 82.95722966666666
$ws.Cells.Item(3, 14).Value = 248.871689
$ws.Cells.Item(3, 15).Value = 0.653306596744776
$ws.Cells.Item(3, 16).Value = 0.653306596744776
$ws.Cells.Item(3, 17).Value = 78.02022370992422
$ws.Cells.Item(3, 18).Value = 702.1820133893179
$ws.Cells.Item(3, 19).Value = 0.01306863897965959
$ws.Cells.Item(3, 20).Value = 0.01306863897965959

# Row 4: ECs -> sCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ncam1"
$ws.Cells.Item(4, 3).Value = "Fgfr1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.9404873333333333
$ws.Cells.Item(4, 8).Value = 2.821462
$ws.Cells.Item(4, 9).Value = 0.02000383747045655
$ws.Cells.Item(4, 10).Value = 0.02000383747045654
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 39.34793966666667
$ws.Cells.Item(4, 14).Value = 118.043819
$ws.Cells.Item(4, 15).Value = 0.3098737585119468
$ws.Cells.Item(4, 16).Value = 0.3098737585119468
$ws.Cells.Item(4, 17).Value = 37.00623884926422
$ws.Cells.Item(4, 18).Value = 333.056149643378
$ws.Cells.Item(4, 19).Value = 0.006198664301632484
$ws.Cells.Item(4, 20).Value = 0.006198664301632482

# Row 5: FAPs -> ECs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ncam1"
$ws.Cells.Item(5, 3).Value = "Fgfr1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.392600333333333
$ws.Cells.Item(5, 8).Value = 4.177801000000001
$ws.Cells.Item(5, 9).Value = 0.0296201232509638
$ws.Cells.Item(5, 10).Value = 0.0296201232509638
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 4.675378666666666
$ws.Cells.Item(5, 14).Value = 14.026136
$ws.Cells.Item(5, 15).Value = 0.03681964474327726
$ws.Cells.Item(5, 16).Value = 0.03681964474327726
$ws.Cells.Item(5, 17).Value = 6.510933889659555
$ws.Cells.Item(5, 18).Value = 58.59840500693601
$ws.Cells.Item(5, 19).Value = 0.001090602415352574
$ws.Cells.Item(5, 20).Value = 0.001090602415352574

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ncam1"
$ws.Cells.Item(6, 3).Value = "Fgfr1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.392600333333333
$ws.Cells.Item(6, 8).Value = 4.177801000000001
$ws.Cells.Item(6, 9).Value = 0.0296201232509638
$ws.Cells.Item(6, 10).Value = 0.0296201232509638
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 82.95722966666666
$ws.Cells.Item(6, 14).Value = 248.871689
$ws.Cells.Item(6, 15).Value = 0.653306596744776
$ws.Cells.Item(6, 16).Value = 0.653306596744776
$ws.Cells.Item(6, 17).Value = 115.5262656862099
$ws.Cells.Item(6, 18).Value = 1039.736391175889
$ws.Cells.Item(6, 19).Value = 0.01935102191624797
$ws.Cells.Item(6, 20).Value = 0.01935102191624797

# Row 7: FAPs -> sCs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ncam1"
$ws.Cells.Item(7, 3).Value = "Fgfr1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.392600333333333
$ws.Cells.Item(7, 8).Value = 4.177801000000001
$ws.Cells.Item(7, 9).Value = 0.0296201232509638
$ws.Cells.Item(7, 10).Value = 0.0296201232509638
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 39.34793966666667
$ws.Cells.Item(7, 14).Value = 118.043819
$ws.Cells.Item(7, 15).Value = 0.3098737585119468
$ws.Cells.Item(7, 16).Value = 0.3098737585119468
$ws.Cells.Item(7, 17).Value = 54.7959538957799
$ws.Cells.Item(7, 18).Value = 493.1635850620191
$ws.Cells.Item(7, 19).Value = 0.009178498919363259
$ws.Cells.Item(7, 20).Value = 0.009178498919363254

# Row 8: sCs -> ECs
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Ncam1"
$ws.Cells.Item(8, 3).Value = "Fgfr1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 44.682258
$ws.Cells.Item(8, 8).Value = 134.046774
$ws.Cells.Item(8, 9).Value = 0.9503760392785797
$ws.Cells.Item(8, 10).Value = 0.9503760392785796
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 4.675378666666666
$ws.Cells.Item(8, 14).Value = 14.026136
$ws.Cells.Item(8, 15).Value = 0.03681964474327726
$ws.Cells.Item(8, 16).Value = 0.03681964474327726
$ws.Cells.Item(8, 17).Value = 208.906475831696
$ws.Cells.Item(8, 18).Value = 1880.158282485264
$ws.Cells.Item(8, 19).Value = 0.03499250813876022
$ws.Cells.Item(8, 20).Value = 0.03499250813876022

# Row 9: sCs -> FAPs
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Ncam1"
$ws.Cells.Item(9, 3).Value = "Fgfr1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 44.682258
$ws.Cells.Item(9, 8).Value = 134.046774
$ws.Cells.Item(9, 9).Value = 0.9503760392785797
$ws.Cells.Item(9, 10).Value = 0.9503760392785796
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 82.95722966666666
$ws.Cells.Item(9, 14).Value = 248.871689
$ws.Cells.Item(9, 15).Value = 0.653306596744776
$ws.Cells.Item(9, 16).Value = 0.653306596744776
$ws.Cells.Item(9, 17).Value = 3706.716338931254
$ws.Cells.Item(9, 18).Value = 33360.44705038128
$ws.Cells.Item(9, 19).Value = 0.6208869358488684
$ws.Cells.Item(9, 20).Value = 0.6208869358488683

# Row 10: sCs -> sCs
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Ncam1"
$ws.Cells.Item(10, 3).Value = "Fgfr1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 44.682258
$ws.Cells.Item(10, 8).Value = 134.046774
$ws.Cells.Item(10, 9).Value = 0.9503760392785797
$ws.Cells.Item(10, 10).Value = 0.9503760392785796
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 39.34793966666667
$ws.Cells.Item(10, 14).Value = 118.043819
$ws.Cells.Item(10, 15).Value = 0.3098737585119468
$ws.Cells.Item(10, 16).Value = 0.3098737585119468
$ws.Cells.Item(10, 17).Value = 1758.154791954434
$ws.Cells.Item(10, 18).Value = 15823.3931275899
$ws.Cells.Item(10, 19).Value = 0.2944965952909511
$ws.Cells.Item(10, 20).Value = 0.294496595290951

